# edit.ps1 - applies the "replace lists" change described by the diff:
#   - {{list.hob1}} {{list.hob2}} -> "Игры Музыка" (run formatting simplified)
#   - {{list.hob3.subhob}} -> "Что-то" (run formatting simplified, bookmark removed)
#   - four new list paragraphs added after that, before the table:
#       "Что-то2" (ilvl 1), "Книги Танцы" (ilvl 0), "Что-то" (ilvl 1), "Что-то2" (ilvl 1)
# The trailing table (with {{table.name1/2/3}}) is intentionally left untouched.

$d = $word.ActiveDocument

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Paragraph 2: {{list.hob1}} {{list.hob2}} -> Игры Музыка ------------------
$p2 = $d.Paragraphs(2)
$p2Xml = $pkgHeader + '<w:p w14:paraId="4D10828C" w14:textId="19C9A04F" w:rsidR="00B460E5" w:rsidRDefault="00B460E5" w:rsidP="00B460E5"><w:pPr><w:pStyle w:val="ae"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="3012"/><w:tab w:val="left" w:pos="3540"/></w:tabs><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Roboto" w:hAnsi="Roboto" w:cs="Roboto"/><w:color w:val="3E3E3E"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="24"/></w:rPr><w:t>Игры Музыка</w:t></w:r></w:p>' + $pkgFooter
$p2.Range.InsertXML($p2Xml) | Out-Null

# --- Paragraph 3: {{list.hob3.subhob}} -> Что-то, bookmark removed ------------
$p3 = $d.Paragraphs(3)
$p3Xml = $pkgHeader + '<w:p w14:paraId="3A1E6A27" w14:textId="559B9680" w:rsidR="00140758" w:rsidRPr="00140758" w:rsidRDefault="00B460E5" w:rsidP="00140758"><w:pPr><w:pStyle w:val="ae"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="8"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="3012"/><w:tab w:val="left" w:pos="3540"/></w:tabs><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Roboto" w:hAnsi="Roboto" w:cs="Roboto"/><w:color w:val="3E3E3E"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:sz w:val="24"/></w:rPr><w:t>Что-то</w:t></w:r></w:p>' + $pkgFooter
$p3.Range.InsertXML($p3Xml) | Out-Null

# --- Insert 4 new list paragraphs right after paragraph 3, before the table ---
function New-ListParaXml([string]$ilvl, [string]$text) {
    $pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
    $pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $pkgHeader + '<w:p><w:pPr><w:pStyle w:val="ae"/><w:numPr><w:ilvl w:val="' + $ilvl + '"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:rPr/><w:t>' + $text + '</w:t></w:r></w:p>' + $pkgFooter
}

$p3 = $d.Paragraphs(3)
$p3.Range.InsertParagraphAfter() | Out-Null

$p4 = $d.Paragraphs(4)
$p4.Range.InsertXML((New-ListParaXml "1" "Что-то2")) | Out-Null
$p4.Range.InsertParagraphAfter() | Out-Null

$p5 = $d.Paragraphs(5)
$p5.Range.InsertXML((New-ListParaXml "0" "Книги Танцы")) | Out-Null
$p5.Range.InsertParagraphAfter() | Out-Null

$p6 = $d.Paragraphs(6)
$p6.Range.InsertXML((New-ListParaXml "1" "Что-то")) | Out-Null
$p6.Range.InsertParagraphAfter() | Out-Null

$p7 = $d.Paragraphs(7)
$p7.Range.InsertXML((New-ListParaXml "1" "Что-то2")) | Out-Null
